# Update init-stock related numbers.
$wb = $excel.ActiveWorkbook

# --- factory sheet: lower the "init count" for the three parts factories ---
$wsFactory = $wb.Worksheets.Item("factory")
$wsFactory.Range("B10").Value = 2
$wsFactory.Range("B11").Value = 3
$wsFactory.Range("B12").Value = 2

# --- stocks sheet: adjust the base stock (column D) feeding the aircraft
#     assembly plant's raw-material rows (17-25, excluding the unchanged
#     17/20/21/24 rows) ---
$wsStocks = $wb.Worksheets.Item("stocks")
$wsStocks.Range("D17").Value = 3200
$wsStocks.Range("D18").Value = 1600
$wsStocks.Range("D19").Value = 800
$wsStocks.Range("D22").Value = 1000
$wsStocks.Range("D23").Value = 4800
$wsStocks.Range("D25").Value = 1200

# --- restore the selections/viewport that Excel records after the edit ---
$wsFactory.Range("B9").Select()
$wsStocks.Activate()
$wsStocks.Application.ActiveWindow.ScrollRow = 13
$wsStocks.Range("D23").Select()
